$wb = $excel.ActiveWorkbook

# --- Content change: "waterfall" sheet, column B (rows 3-16) ---
# Was referencing the shared string "归集账户" (collection account),
# now referencing "账户I" (account I).
$wsWaterfall = $wb.Worksheets.Item("waterfall")
$wsWaterfall.Range("B3:B16").Value = "账户I"

# --- View/selection changes ---
# "account" sheet used to be the active/selected tab with selection C4;
# it is no longer the selected tab and its selection moves to A3.
$wsAccount = $wb.Worksheets.Item("account")
[void]$wsAccount.Range("A3").Select()

# "waterfall" sheet becomes the active/selected tab, with selection C6
# (previously selection was J16, and it was not the active tab).
[void]$wsWaterfall.Activate()
[void]$wsWaterfall.Range("C6").Select()
